# Quarterly indexing esoteric bug-fix operation
# The forecast-error table was re-indexed by one quarter: every row's
# error metrics (ME, MAE, MSE, RMSE, SE) and observation count (N) are
# recomputed against the corrected quarterly offset, while the Q-label
# column (A) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1656278310671722
$ws.Range("C2").Value = 0.4719345280589418
$ws.Range("D2").Value = 0.5088034263367071
$ws.Range("E2").Value = 0.7133045817438067
$ws.Range("F2").Value = 0.7022190543561216
$ws.Range("G2").Value = 42

$ws.Range("B3").Value = -0.01237757133353611
$ws.Range("C3").Value = 0.5576922206168767
$ws.Range("D3").Value = 0.4968421630465693
$ws.Range("E3").Value = 0.7048703164742924
$ws.Range("F3").Value = 0.7135167711720684
$ws.Range("G3").Value = 41

$ws.Range("B4").Value = 0.182806223233343
$ws.Range("C4").Value = 0.5882942346078044
$ws.Range("D4").Value = 0.5943299325354084
$ws.Range("E4").Value = 0.7709279684480311
$ws.Range("F4").Value = 0.7584814905927938
$ws.Range("G4").Value = 40

$ws.Range("B5").Value = 0.04193216340107105
$ws.Range("C5").Value = 0.6060728462292108
$ws.Range("D5").Value = 0.5212737734367096
$ws.Range("E5").Value = 0.7219929178577236
$ws.Range("F5").Value = 0.7301964987385136
$ws.Range("G5").Value = 39

$ws.Range("B6").Value = 0.1893694583833012
$ws.Range("C6").Value = 0.6581235952936254
$ws.Range("D6").Value = 0.6823089728457381
$ws.Range("E6").Value = 0.8260199591085787
$ws.Range("F6").Value = 0.8148127107126543
$ws.Range("G6").Value = 38

$ws.Range("B7").Value = 0.09660088401069185
$ws.Range("C7").Value = 0.5163580642747845
$ws.Range("D7").Value = 0.4277032691171462
$ws.Range("E7").Value = 0.6539902668367062
$ws.Range("F7").Value = 0.6557384920421035
$ws.Range("G7").Value = 37

$ws.Range("B8").Value = 0.1755910946091918
$ws.Range("C8").Value = 0.6654140680913321
$ws.Range("D8").Value = 0.6931489096200038
$ws.Range("E8").Value = 0.8325556495634413
$ws.Range("F8").Value = 0.8253726496230516
$ws.Range("G8").Value = 36

$ws.Range("B9").Value = 0.1322577806459256
$ws.Range("C9").Value = 0.6219002833064226
$ws.Range("D9").Value = 0.6041857322815641
$ws.Range("E9").Value = 0.7772938519514767
$ws.Range("F9").Value = 0.7771417542528034
$ws.Range("G9").Value = 35

$ws.Range("B10").Value = 0.2251170054142814
$ws.Range("C10").Value = 0.7236725851673381
$ws.Range("D10").Value = 0.8349781460265824
$ws.Range("E10").Value = 0.9137713860843872
$ws.Range("F10").Value = 0.8989255592700524
$ws.Range("G10").Value = 34

$ws.Range("B11").Value = 0.1090013000149534
$ws.Range("C11").Value = 0.6144815156716034
$ws.Range("D11").Value = 0.5898028043825534
$ws.Range("E11").Value = 0.7679862006459187
$ws.Range("F11").Value = 0.7719984251979752
$ws.Range("G11").Value = 33
